# Update "Hjemme passive" data: recomputed meanEMG / legmaxROM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header/sample-size row) B1:E1
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON)
$ws.Range("B2").Value = 8.2258820116076379
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 10.872399275861497
$ws.Range("E2").ClearContents()

# Row 3 (STR)
$ws.Range("B3").Value = 7.5777447545056162
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 9.6826418635996561
$ws.Range("E3").Value = -7.3463355315719454

# Update the selected range shown in the sheet view
$ws.Range("B1:E3").Select()
